$d = $word.ActiveDocument

# Update the date paragraph (first paragraph in the document)
$d.Paragraphs.Item(1).Range.Text = "2025-06-22 Sunday"

# Update each math-equation cell in the 20x5 table, in row-major order
$t = $d.Tables.Item(1)
$values = @(
    "12+59=",
    "3+38=",
    "68+7=",
    "35+29=",
    "9+39=",
    "16-7=",
    "84-37=",
    "17+29=",
    "55+29=",
    "91-88=",
    "33-29=",
    "66+27=",
    "23+69=",
    "32-16=",
    "61-27=",
    "56+5=",
    "84-55=",
    "24-7=",
    "48+39=",
    "81-65=",
    "61-44=",
    "69+24=",
    "44+8=",
    "35-9=",
    "94-39=",
    "29+4=",
    "91-78=",
    "80-42=",
    "90-61=",
    "66-49=",
    "64-57=",
    "97-18=",
    "47+5=",
    "15+58=",
    "74-49=",
    "66-48=",
    "85-66=",
    "59+15=",
    "69+4=",
    "25+57=",
    "50-16=",
    "92-9=",
    "11-2=",
    "69+6=",
    "60-13=",
    "9+83=",
    "54+9=",
    "32-4=",
    "91-66=",
    "19+7=",
    "29+54=",
    "3+78=",
    "91-19=",
    "85-57=",
    "81-34=",
    "15+48=",
    "98-19=",
    "59+18=",
    "29+38=",
    "55-7=",
    "12+19=",
    "69+14=",
    "75+9=",
    "74-59=",
    "35+27=",
    "60-36=",
    "73-6=",
    "81-18=",
    "25+46=",
    "70-16=",
    "85-18=",
    "98-9=",
    "40-39=",
    "77+15=",
    "5+89=",
    "38-19=",
    "63-29=",
    "71-9=",
    "45-7=",
    "81-29=",
    "48+28=",
    "9+26=",
    "8+43=",
    "51-16=",
    "54-39=",
    "23+28=",
    "34+28=",
    "8+19=",
    "55+16=",
    "66+16=",
    "72-27=",
    "13+59=",
    "77-28=",
    "78+9=",
    "20-12=",
    "72-64=",
    "70-51=",
    "9+39=",
    "46+35=",
    "55-39="
)

$rows = 20
$cols = 5
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Done. Updated" $idx "cells."